$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "PARTNER & STRATEGIC CONSULTANT - Siege Analytics, Washington, DC | January 2014 – Present"; New = "PARTNER - Siege Analytics, Washington, DC | January 2014 – Present" },
    @{ Old = "PRINCIPAL MARKETING CONSULTANT - Clarity and Rigour, Washington, DC | 2012 – 2014"; New = "DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 – 2014" },
    @{ Old = "DIRECTOR OF MARKETING - Helm, Washington, DC | 2010 – 2012"; New = "SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 – 2012" },
    @{ Old = "SENIOR MARKETING ANALYST - GSD&M, Austin, TX | 2008 – 2010"; New = "SENIOR ANALYST - Myers Research, Washington, DC | 2008 – 2010" },
    @{ Old = "MARKETING COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008"; New = "RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008" },
    @{ Old = "MARKETING SPECIALIST - Salsa Labs, Inc., Washington, DC | 2004 – 2006"; New = "SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 – 2006" },
    @{ Old = "COMMUNICATIONS COORDINATOR - The Praxis Project, Oakland, CA | 2002 – 2004"; New = "INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004" },
    @{ Old = "RESEARCH COORDINATOR - Lake Research Partners, Washington, DC | 2001 – 2002"; New = "PROGRAMMER - Lake Research Partners, Washington, DC | 2001 – 2002" },
    @{ Old = "FIELD COORDINATOR - The Feldman Group, Washington, DC | 2000 – 2001"; New = "FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 – 2001" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}

Write-Host "Done applying job title replacements."
